$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row Price (D) / Volume(1h) (E) updates.
# Price values are assigned with a leading apostrophe so Excel keeps them as
# literal text (matching the original inlineStr cells) instead of coercing
# numeric-looking strings like "238.71" into real numbers; Style is reset to
# 'Normal' right after so no stray number-format style sticks to the cell.
$ws.Range("D2").Value = "'" + '41.614.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").Value = "'" + '2.168.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.11%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'" + '238.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("D6").Value = "'" + '0.608'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.38%  '
$ws.Range("D7").Value = "'" + '72.27'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.24%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = "'" + '0.577'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.55%  '
$ws.Range("D10").Value = "'" + '39.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.13%  '
$ws.Range("D11").Value = "'" + '0.0908'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.06%  '
$ws.Range("D12").Value = "'" + '54.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.37%  '
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("D14").Value = "'" + '6.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.29%  '
$ws.Range("D15").Value = "'" + '2.493.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("D16").Value = "'" + '14.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.43%  '
$ws.Range("D17").Value = "'" + '2.153.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.73%  '
$ws.Range("D18").Value = "'" + '0.782'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.29%  '
$ws.Range("D19").Value = "'" + '41.498.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("D20").Value = "'" + '0.0000104'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.56%  '
$ws.Range("D21").Value = "'" + '69.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.49%  '
$ws.Range("D22").Value = "'" + '5.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.61%  '
$ws.Range("D23").Value = "'" + '9.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -13.47%  '
$ws.Range("D24").Value = "'" + '227.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'" + '2.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.36%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = "'" + '10.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.66%  '
$ws.Range("E28").Value = '  -9.88%  '
$ws.Range("D29").Value = "'" + '2.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.26%  '
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("D31").Value = "'" + '170.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("D32").Value = "'" + '19.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.21%  '
$ws.Range("D33").Value = "'" + '33.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.55%  '
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("E35").Value = '  -10.73%  '
$ws.Range("E36").Value = '  -3.62%  '
$ws.Range("D37").Value = "'" + '4.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("E38").Value = '  -4.95%  '
$ws.Range("D39").Value = "'" + '0.0304'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("E40").Value = '  -2.13%  '
$ws.Range("D41").Value = "'" + '11.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -12.51%  '
$ws.Range("D42").Value = "'" + '5.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.41%  '
$ws.Range("D43").Value = "'" + '59.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.44%  '
$ws.Range("D46").Value = "'" + '0.0961'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.45%  '
$ws.Range("D47").Value = "'" + '96.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.83%  '
$ws.Range("D48").Value = "'" + '1.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("E49").Value = '  -5.73%  '
$ws.Range("E50").Value = '  -8.16%  '

# Rows 44 and 45 swap places (Algorand <-> FraxShare) with updated values
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = "'" + '8.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.98%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = "'" + '0.189'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.86%  '
